# #3456 updated PM Property ID
# Updates the "Portfolio Manager Building ID" values in column B (rows 2-10)
# of the BPS Data sheet, and restores the sheet selection to B2:B10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPropertyIds = @(22178843, 22178844, 22178845, 22178846, 22178847, 22178848, 22178849, 22178850, 22178851)

$row = 2
foreach ($id in $newPropertyIds) {
    $ws.Cells.Item($row, 2).Value = $id
    $row++
}

# Match the workbook's saved selection state (B2:B10, active cell B2)
$ws.Range("B2:B10").Select()
